# chgt xxx colonne export + selection ihm
#
# The "Profil" header label (B1) gets a temporary/debug suffix "XXXXX",
# and the active selection on the sheet moves from D14 to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Profil" column header to "ProfilXXXXX"
$ws.Range("B1").Value = "ProfilXXXXX"

# Update the current selection shown in the UI
$ws.Range("D12").Select()
